$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "SCD0023"

# Update TC_ID value in B2
$ws.Range("B2").Value = "SCD0023-001"

# Re-fit column B width to its (now longer) content (Excel auto bestFit
# recalculated this when the B2 text got longer)
$ws.Columns.Item(2).ColumnWidth = 11.6

# Select B3 (updates active cell / selection, also resets scroll position)
$ws.Range("B3").Select()
